$wb = $excel.ActiveWorkbook

# --- Update data on "Proposed PVs (NoRaid - Raw)" sheet ---
$wsRaw = $wb.Worksheets.Item("Proposed PVs (NoRaid - Raw)")

# Re-RAID disks: raw disk sizes go from 3TB to 3.5TB for rows 5-8 (B:G)
$wsRaw.Range("B5:G8").Value = 3.5

# Re-enter the summary formulas as a single multi-cell assignment so Excel
# stores them as shared formulas (matching a real Excel re-save/fill).
$wsRaw.Range("B9:G9").Formula = "=SUM(B5:B8)"
$wsRaw.Range("B10:H10").Formula = "=B3-B9"

# --- Update sheet selections / active tab ---
# Select A5 on the "NoRaid - Raw" sheet (it will no longer be the active tab)
$wsRaw.Activate()
$wsRaw.Range("A5").Select()

# Make "Proposed PVs" the active sheet/tab, keeping its existing B5 selection
$wsProposed = $wb.Worksheets.Item("Proposed PVs")
$wsProposed.Activate()
$wsProposed.Range("B5").Select()
